$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 12:22"

# Austria (row 23)
$ws.Range("B23").Value = 14795
$ws.Range("C23").Value = 46
$ws.Range("D23").Value = 10631
$ws.Range("E23").Value = 3694
$ws.Range("F23").Value = 194
$ws.Range("G23").Value = 18
$ws.Range("H23").Value = 470

# Rumania (row 32)
$ws.Range("B32").Value = 8936
$ws.Range("C32").Value = 190
$ws.Range("D32").Value = 2017
$ws.Range("E32").Value = 6468
$ws.Range("F32").Value = 261

# Uzbekistan (row 68)
$ws.Range("B68").Value = 1582
$ws.Range("C68").Value = 17
$ws.Range("D68").Value = 238
$ws.Range("E68").Value = 1339

# Hong Kong (row 83)
$ws.Range("D83").Value = 630
$ws.Range("E83").Value = 392
